{"js": "// Update the 25 two-digit-division answer cells in the practice table.\n// Each [oldText, newText] pair is keyed to exactly one unique <w:t> run,\n// so a document-wide exact-text search unambiguously targets the right cell.\nconst replacements = [\n  [\"78\u00f74=19, 2\", \"18\u00f76=3, 0\"],\n  [\"52\u00f74=13, 0\", \"56\u00f79=6, 2\"],\n  [\"14\u00f78=1, 6\", \"12\u00f77=1, 5\"],\n  [\"83\u00f76=13, 5\", \"13\u00f78=1, 5\"],\n  [\"71\u00f76=11, 5\", \"21\u00f73=7, 0\"],\n  [\"19\u00f76=3, 1\", \"81\u00f73=27, 0\"],\n  [\"39\u00f72=19, 1\", \"27\u00f72=13, 1\"],\n  [\"55\u00f78=6, 7\", \"23\u00f73=7, 2\"],\n  [\"81\u00f79=9, 0\", \"58\u00f73=19, 1\"],\n  [\"22\u00f79=2, 4\", \"26\u00f75=5, 1\"],\n  [\"45\u00f77=6, 3\", \"26\u00f74=6, 2\"],\n  [\"69\u00f79=7, 6\", \"43\u00f77=6, 1\"],\n  [\"89\u00f78=11, 1\", \"12\u00f73=4, 0\"],\n  [\"50\u00f76=8, 2\", \"65\u00f76=10, 5\"],\n  [\"21\u00f78=2, 5\", \"34\u00f74=8, 2\"],\n  [\"84\u00f75=16, 4\", \"36\u00f77=5, 1\"],\n  [\"28\u00f78=3, 4\", \"90\u00f74=22, 2\"],\n  [\"47\u00f73=15, 2\", \"44\u00f75=8, 4\"],\n  [\"20\u00f78=2, 4\", \"57\u00f74=14, 1\"],\n  [\"27\u00f76=4, 3\", \"33\u00f78=4, 1\"],\n  [\"85\u00f76=14, 1\", \"79\u00f75=15, 4\"],\n  [\"15\u00f78=1, 7\", \"75\u00f72=37, 1\"],\n  [\"63\u00f73=21, 0\", \"20\u00f79=2, 2\"],\n  [\"38\u00f77=5, 3\", \"40\u00f79=4, 4\"],\n  [\"61\u00f78=7, 5\", \"42\u00f72=21, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit-division answer cells with their new values.\n# Each (old, new) pair corresponds 1:1 to a unique <w:t> run in the table,\n# so a document-wide Find/Replace targeting the exact old text is unambiguous.\n$pairs = @(\n    @(\"78\u00f74=19, 2\", \"18\u00f76=3, 0\"),\n    @(\"52\u00f74=13, 0\", \"56\u00f79=6, 2\"),\n    @(\"14\u00f78=1, 6\", \"12\u00f77=1, 5\"),\n    @(\"83\u00f76=13, 5\", \"13\u00f78=1, 5\"),\n    @(\"71\u00f76=11, 5\", \"21\u00f73=7, 0\"),\n    @(\"19\u00f76=3, 1\", \"81\u00f73=27, 0\"),\n    @(\"39\u00f72=19, 1\", \"27\u00f72=13, 1\"),\n    @(\"55\u00f78=6, 7\", \"23\u00f73=7, 2\"),\n    @(\"81\u00f79=9, 0\", \"58\u00f73=19, 1\"),\n    @(\"22\u00f79=2, 4\", \"26\u00f75=5, 1\"),\n    @(\"45\u00f77=6, 3\", \"26\u00f74=6, 2\"),\n    @(\"69\u00f79=7, 6\", \"43\u00f77=6, 1\"),\n    @(\"89\u00f78=11, 1\", \"12\u00f73=4, 0\"),\n    @(\"50\u00f76=8, 2\", \"65\u00f76=10, 5\"),\n    @(\"21\u00f78=2, 5\", \"34\u00f74=8, 2\"),\n    @(\"84\u00f75=16, 4\", \"36\u00f77=5, 1\"),\n    @(\"28\u00f78=3, 4\", \"90\u00f74=22, 2\"),\n    @(\"47\u00f73=15, 2\", \"44\u00f75=8, 4\"),\n    @(\"20\u00f78=2, 4\", \"57\u00f74=14, 1\"),\n    @(\"27\u00f76=4, 3\", \"33\u00f78=4, 1\"),\n    @(\"85\u00f76=14, 1\", \"79\u00f75=15, 4\"),\n    @(\"15\u00f78=1, 7\", \"75\u00f72=37, 1\"),\n    @(\"63\u00f73=21, 0\", \"20\u00f79=2, 2\"),\n    @(\"38\u00f77=5, 3\", \"40\u00f79=4, 4\"),\n    @(\"61\u00f78=7, 5\", \"42\u00f72=21, 0\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
